$wb = $excel.ActiveWorkbook

$newShortName = "247-MS-EI-DB-SAR-REC-NON-RNI-CTRFD-SAR-MD-TR-1-ONTIME"

$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Update the product short name value on both sheets (B1)
$wsInput.Range("B1").Value = $newShortName
$wsOutput.Range("B1").Value = $newShortName

# Move the active cell selection to B1 on both sheets
$wsInput.Range("B1").Select()
$wsOutput.Range("B1").Select()

# Make ProductLoanOutput the active/selected sheet (and tab)
$wsOutput.Activate()
